$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: new section header "Memory Mapping (40 Bytes)" with the same style as the
# other section headers in row 1 (merged A:C, yellow highlighted, bold+underline font)
$ws.Range("A1:C1").Copy()
$ws.Range("A14:C14").PasteSpecial(-4122)
$ws.Range("A14").Value = "Memory Mapping (40 Bytes)"
$ws.Range("A14:C14").Merge()
$ws.Rows.Item(14).RowHeight = 18.75

# Row 15: column headers for the new table, matching the existing "Run " / "Time (ns)" headers
$ws.Range("A15").Value = "Run "
$ws.Range("B15").Value = "Time (ns)"

# Rows 16-25: run numbers 1 through 10
for ($i = 1; $i -le 10; $i++) {
    $ws.Cells.Item(15 + $i, 1).Value = $i
}

# Row 26: "Average" label (value not yet filled in)
$ws.Range("A26").Value = "Average"

# Leave the selection where the data entry stopped
[void]$ws.Range("E19").Select()
